# The event "合肥·排球少年only之夏日招新季" (row 2) has ended / been removed from
# the feed. Remove its row from every sheet that lists it ("展览" and
# "全部类型"); the remaining rows shift up to fill the gap. A handful of
# "想去人数" (interest count) values were also refreshed upstream, and the
# "合肥·《四月是你的谎言》" listing's count ticked up by one (visible in both
# the "演出" sheet and the "全部类型" sheet).

$wb = $excel.ActiveWorkbook

function Remove-EventRow {
    param($ws)

    # Delete the whole row - remaining rows (and their data) shift up by one.
    $ws.Rows(2).Delete()

    # Column A is just a sequential 0-based row index (header=0, first data
    # row=1, ...) - not event-specific data - so re-number it after the
    # shift instead of leaving the gap left behind by the deleted row.
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

# --- 展览 (sheet 1) ---
$wsExpo = $wb.Worksheets.Item(1)
Remove-EventRow $wsExpo
$wsExpo.Range("F2").Value = 118     # 合肥·次元盛典燃动DNA
$wsExpo.Range("F4").Value = 5101    # 合肥·第八届环形宇宙动漫游戏嘉年华Plus
$wsExpo.Range("F8").Value = 768     # 合肥·SSS第五人格only
$wsExpo.Range("F9").Value = 252     # 合肥·国乙only宇宙心动（含夜场）

# --- 演出 (sheet 2) --- no rows removed, just the refreshed interest count
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F3").Value = 5       # 合肥·《四月是你的谎言》…

# --- 本地生活 (sheet 3) --- unchanged

# --- 全部类型 (sheet 4) ---
$wsAll = $wb.Worksheets.Item(4)
Remove-EventRow $wsAll
$wsAll.Range("F2").Value = 118      # 合肥·次元盛典燃动DNA
$wsAll.Range("F4").Value = 5101     # 合肥·第八届环形宇宙动漫游戏嘉年华Plus
$wsAll.Range("F8").Value = 768      # 合肥·SSS第五人格only
$wsAll.Range("F10").Value = 252     # 合肥·国乙only宇宙心动（含夜场）
$wsAll.Range("F12").Value = 5       # 合肥·《四月是你的谎言》…
